$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ex1_data")
$ws1.Range("A1").Value = "test"
